# Auto-generated Excel COM-interop edit script
# Updates cryptos list data per target diff (Wed Nov 15 23:35:23 UTC 2023)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text / non-numeric-looking values: direct assignment ---
$ws.Range("D2").Value = "37.770.45"
$ws.Range("E2").Value = "  +5.83%  "
$ws.Range("D3").Value = "2.051.55"
$ws.Range("E3").Value = "  +3.31%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("E5").Value = "  +3.79%  "
$ws.Range("E6").Value = "  +2.13%  "
$ws.Range("E7").Value = "  +13.55%  "
$ws.Range("E8").Value = "  -0.12%  "
$ws.Range("E9").Value = "  +2.00%  "
$ws.Range("E10").Value = "  +4.83%  "
$ws.Range("E11").Value = "  +4.38%  "
$ws.Range("E12").Value = "  +1.00%  "
$ws.Range("E13").Value = "  -0.85%  "
$ws.Range("E14").Value = "  +7.41%  "
$ws.Range("D15").Value = "2.349.89"
$ws.Range("E15").Value = "  +3.19%  "
$ws.Range("E16").Value = "  +20.12%  "
$ws.Range("E17").Value = "  +5.62%  "
$ws.Range("D18").Value = "2.045.31"
$ws.Range("E18").Value = "  +3.00%  "
$ws.Range("D19").Value = "37.656.74"
$ws.Range("E19").Value = "  +5.89%  "
$ws.Range("E20").Value = "  +4.96%  "
$ws.Range("D21").Value = "0.0₃0877"
$ws.Range("E21").Value = "  +4.50%  "
$ws.Range("E22").Value = "  +5.59%  "
$ws.Range("E23").Value = "  +2.33%  "
$ws.Range("E24").Value = "  +14.34%  "
$ws.Range("E25").Value = "  +0.10%  "
$ws.Range("E26").Value = "  +4.55%  "
$ws.Range("E27").Value = "  +5.17%  "
$ws.Range("E28").Value = "  -2.23%  "
$ws.Range("E29").Value = "  +2.45%  "
$ws.Range("B30").Value = "Kaspa"
$ws.Range("C30").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("E30").Value = "  +27.00%  "
$ws.Range("B31").Value = "Stellar"
$ws.Range("C31").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("E31").Value = "  +2.42%  "
$ws.Range("E32").Value = "  +8.61%  "
$ws.Range("E33").Value = "  +6.56%  "
$ws.Range("E34").Value = "  +10.88%  "
$ws.Range("E35").Value = "  +4.85%  "
$ws.Range("E36").Value = "  +2.44%  "
$ws.Range("E37").Value = "  +3.26%  "
$ws.Range("E38").Value = "  -0.09%  "
$ws.Range("E39").Value = "  +23.25%  "
$ws.Range("E40").Value = "  +16.66%  "
$ws.Range("E41").Value = "  +24.95%  "
$ws.Range("E42").Value = "  +3.79%  "
$ws.Range("E43").Value = "  +4.33%  "
$ws.Range("E44").Value = "  +3.73%  "
$ws.Range("E45").Value = "  +5.34%  "
$ws.Range("E46").Value = "  +9.85%  "
$ws.Range("E47").Value = "  +7.50%  "
$ws.Range("E48").Value = "  +4.48%  "
$ws.Range("D49").Value = "1.411.28"
$ws.Range("E49").Value = "  +2.35%  "
$ws.Range("E50").Value = "  +2.18%  "
$ws.Range("E51").Value = "  +3.15%  "

# --- Numeric-looking values that must remain text: force text format, then restore default style ---
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "252.08"
$cell.Style = "Normal"
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "0.653"
$cell.Style = "Normal"
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "65.20"
$cell.Style = "Normal"
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.999"
$cell.Style = "Normal"
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "60.98"
$cell.Style = "Normal"
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.0764"
$cell.Style = "Normal"
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "0.919"
$cell.Style = "Normal"
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "15.15"
$cell.Style = "Normal"
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "74.18"
$cell.Style = "Normal"
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "238.72"
$cell.Style = "Normal"
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "2.69"
$cell.Style = "Normal"
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "2.40"
$cell.Style = "Normal"
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "9.63"
$cell.Style = "Normal"
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "159.98"
$cell.Style = "Normal"
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "19.98"
$cell.Style = "Normal"
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "0.114"
$cell.Style = "Normal"
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "0.122"
$cell.Style = "Normal"
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "5.22"
$cell.Style = "Normal"
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "4.73"
$cell.Style = "Normal"
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "0.0619"
$cell.Style = "Normal"
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "1.87"
$cell.Style = "Normal"
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "0.0220"
$cell.Style = "Normal"
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "17.00"
$cell.Style = "Normal"
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "8.01"
$cell.Style = "Normal"
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "95.19"
$cell.Style = "Normal"
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "2.95"
$cell.Style = "Normal"
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "47.29"
$cell.Style = "Normal"
